$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 572.3333
$ws.Range("J38").Value = 4000
$ws.Range("L38").Value = 12000
$ws.Range("N38").Value = -12744
$ws.Range("H58").Value = 4412.0713
$ws.Range("J58").Value = 8507
$ws.Range("L58").Value = 25521
$ws.Range("N58").Value = -25821
$ws.Range("H106").Value = 6957.2
$ws.Range("I106").Value = 4445.75
$ws.Range("K106").Value = 4445.75
$ws.Range("M106").Value = -3814.75
$ws.Range("H137").Value = 890562.2
$ws.Range("I137").Value = 1357.6111
$ws.Range("J137").Value = 3177088.2
$ws.Range("K137").Value = 4072.8333
$ws.Range("L137").Value = 9531264.600000001
$ws.Range("M137").Value = -1522.8333
$ws.Range("N137").Value = -9536364.600000001
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 852.43335
$ws.Range("I2").Value = 852.43335
$ws.Range("K2").Value = 852.43335
$ws.Range("M2").Value = -739.43335
$ws.Range("H32").Value = 16926.895
$ws.Range("I32").Value = 19483.074
$ws.Range("J32").Value = 5424.0835
$ws.Range("K32").Value = 19483.074
$ws.Range("L32").Value = 5424.0835
$ws.Range("M32").Value = -19196.074
$ws.Range("N32").Value = -5998.0835
$ws.Range("H45").Value = 2662.3333
$ws.Range("I45").Value = 1763.9
$ws.Range("J45").Value = 3785.375
$ws.Range("K45").Value = 1763.9
$ws.Range("L45").Value = 3785.375
$ws.Range("M45").Value = -1386.9
$ws.Range("N45").Value = -4539.375
$ws.Range("H61").Value = 2726.2354
$ws.Range("I61").Value = 2438
$ws.Range("K61").Value = 2438
$ws.Range("M61").Value = -2226
$ws.Range("H96").Value = 37364.57
$ws.Range("J96").Value = 37364.57
$ws.Range("L96").Value = 37364.57
$ws.Range("N96").Value = -42856.57
$ws.Range("H116").Value = 852.43335
$ws.Range("I116").Value = 852.43335
$ws.Range("K116").Value = 852.43335
$ws.Range("M116").Value = 1441.56665
$ws.Range("H122").Value = 1820.85
$ws.Range("I122").Value = 1502.0714
$ws.Range("J122").Value = 2564.6667
$ws.Range("K122").Value = 4506.2142
$ws.Range("L122").Value = 7694.000100000001
$ws.Range("M122").Value = -2056.2142
$ws.Range("N122").Value = -12594.0001
$ws.Range("H132").Value = 26566.117
$ws.Range("I132").Value = 33593.28
$ws.Range("K132").Value = 100779.84
$ws.Range("M132").Value = -98249.84
$ws.Range("H136").Value = 2726.2354
$ws.Range("I136").Value = 2438
$ws.Range("K136").Value = 7314
$ws.Range("M136").Value = -4764
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 852.43335
$ws.Range("I3").Value = 852.43335
$ws.Range("K3").Value = 852.43335
$ws.Range("M3").Value = -738.43335
$ws.Range("H105").Value = 4666.5
$ws.Range("I105").Value = 4500
$ws.Range("K105").Value = 4500
$ws.Range("M105").Value = -2753
$ws.Range("H125").Value = 92000
$ws.Range("J125").Value = 92000
$ws.Range("L125").Value = 92000
$ws.Range("N125").Value = -101840
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 224.83333
$ws.Range("I7").Value = 279.875
$ws.Range("J7").Value = 114.75
$ws.Range("K7").Value = 279.875
$ws.Range("L7").Value = 114.75
$ws.Range("M7").Value = -166.875
$ws.Range("N7").Value = -340.75
$ws.Range("H16").Value = 3157.4443
$ws.Range("I16").Value = 2733.5833
$ws.Range("K16").Value = 2733.5833
$ws.Range("M16").Value = -2446.5833
$ws.Range("H58").Value = 102889.1
$ws.Range("I58").Value = 113932.445
$ws.Range("K58").Value = 113932.445
$ws.Range("M58").Value = -113729.445
$ws.Range("H107").Value = 2050.9375
$ws.Range("I107").Value = 201.66667
$ws.Range("K107").Value = 201.66667
$ws.Range("M107").Value = 1718.33333
$ws.Range("H113").Value = 3157.4443
$ws.Range("I113").Value = 2733.5833
$ws.Range("K113").Value = 2733.5833
$ws.Range("M113").Value = -563.5832999999998
$ws.Range("H134").Value = 79065.53999999999
$ws.Range("I134").Value = 85537.586
$ws.Range("K134").Value = 256612.758
$ws.Range("M134").Value = -254077.758
$ws.Range("H136").Value = 102889.1
$ws.Range("I136").Value = 113932.445
$ws.Range("K136").Value = 341797.335
$ws.Range("M136").Value = -339247.335
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 477261.53
$ws.Range("J107").Value = 527441.7
$ws.Range("L107").Value = 1582325.1
$ws.Range("N107").Value = -1586165.1
$ws.Range("H113").Value = 1617.7307
$ws.Range("J113").Value = 4779.6
$ws.Range("L113").Value = 14338.8
$ws.Range("N113").Value = -18678.8
$ws.Range("H136").Value = 5663
$ws.Range("I136").Value = 3494.5
$ws.Range("K136").Value = 10483.5
$ws.Range("M136").Value = -5383.5
$ws.Range("H139").Value = 6631.1177
$ws.Range("I139").Value = 946.6
$ws.Range("K139").Value = 2839.8
$ws.Range("M139").Value = 2300.2
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 22847
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()
$ws.Range("H80").Value = 7182.7646
$ws.Range("I80").Value = 5687.091
$ws.Range("J80").Value = 9924.833000000001
$ws.Range("K80").Value = 5687.091
$ws.Range("L80").Value = 9924.833000000001
$ws.Range("M80").Value = -4689.091
$ws.Range("N80").Value = -11920.833
$ws.Range("H83").Value = 7182.7646
$ws.Range("I83").Value = 5687.091
$ws.Range("J83").Value = 9924.833000000001
$ws.Range("K83").Value = 28435.455
$ws.Range("L83").Value = 49624.165
$ws.Range("M83").Value = -23443.455
$ws.Range("N83").Value = -59608.165
$ws.Range("H102").Value = 2282.7144
$ws.Range("I102").Value = 2498.1667
$ws.Range("K102").Value = 2498.1667
$ws.Range("M102").Value = -876.1667000000002
$ws.Range("H113").Value = 4257.8
$ws.Range("I113").Value = 5644.5
$ws.Range("J113").Value = 3333.3333
$ws.Range("K113").Value = 5644.5
$ws.Range("L113").Value = 3333.3333
$ws.Range("M113").Value = -3474.5
$ws.Range("N113").Value = -7673.3333
$ws.Range("H122").Value = 3335.32
$ws.Range("I122").Value = 1969.125
$ws.Range("K122").Value = 5907.375
$ws.Range("M122").Value = -3457.375
$ws.Range("H132").Value = 41939.69
$ws.Range("J132").Value = 2646.3333
$ws.Range("L132").Value = 7938.999899999999
$ws.Range("N132").Value = -12998.9999
$ws.Range("H133").Value = 145852.33
$ws.Range("J133").Value = 145852.33
$ws.Range("L133").Value = 145852.33
$ws.Range("N133").Value = -155972.33
$ws.Range("H134").Value = 44999.25
$ws.Range("J134").Value = 44999.25
$ws.Range("L134").Value = 134997.75
$ws.Range("N134").Value = -140067.75
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("H55").Value = 394
$ws.Range("I55").Value = 306.54544
$ws.Range("K55").Value = 306.54544
$ws.Range("M55").Value = -133.54544
$ws.Range("H61").Value = 6558.975
$ws.Range("I61").Value = 5813.6772
$ws.Range("J61").Value = 9126.111000000001
$ws.Range("K61").Value = 5813.6772
$ws.Range("L61").Value = 9126.111000000001
$ws.Range("M61").Value = -5611.6772
$ws.Range("N61").Value = -9530.111000000001
$ws.Range("H100").Value = 7624.9165
$ws.Range("J100").Value = 7758.4287
$ws.Range("L100").Value = 7758.4287
$ws.Range("N100").Value = -8840.4287
$ws.Range("H113").Value = 6558.975
$ws.Range("I113").Value = 5813.6772
$ws.Range("J113").Value = 9126.111000000001
$ws.Range("K113").Value = 5813.6772
$ws.Range("L113").Value = 9126.111000000001
$ws.Range("M113").Value = -3643.6772
$ws.Range("N113").Value = -13466.111
$ws.Range("H136").Value = 5125.1665
$ws.Range("J136").Value = 7971.25
$ws.Range("L136").Value = 23913.75
$ws.Range("N136").Value = -29013.75
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 710.5789
$ws.Range("I100").Value = 620.2
$ws.Range("K100").Value = 1240.4
$ws.Range("M100").Value = -699.4000000000001
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H122").Value = 4296.4614
$ws.Range("I122").Value = 4839.5
$ws.Range("J122").Value = 2486.3333
$ws.Range("K122").Value = 14518.5
$ws.Range("L122").Value = 7458.999899999999
$ws.Range("M122").Value = -12068.5
$ws.Range("N122").Value = -12358.9999
$ws.Range("H126").Value = 4902.0884
$ws.Range("I126").Value = 4734.2964
$ws.Range("K126").Value = 14202.8892
$ws.Range("M126").Value = -11732.8892
$ws.Range("H132").Value = 29109.783
$ws.Range("I132").Value = 29723.945
$ws.Range("K132").Value = 89171.83499999999
$ws.Range("M132").Value = -86641.83499999999
$ws.Range("H136").Value = 3945.8572
$ws.Range("I136").Value = 2814.818
$ws.Range("K136").Value = 8444.454000000002
$ws.Range("M136").Value = -5894.454000000002
$ws.Range("H138").Value = 149999
$ws.Range("I138").Value = 149999
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 149999
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = -144859
$ws.Range("N138").ClearContents()
